$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 97, shifting the existing data (rows 97:152) down by one
# (new last row becomes 153), matching the weekly update to the price series.
$ws.Rows(97).Insert()

# Populate the newly inserted row 97 with the new weekly data point.
$ws.Range("A97").Value = 4
$ws.Range("B97").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C97").Value = "Los Lagos"
$ws.Range("D97").Value = 44438
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = 100112037
$ws.Range("G97").Value = "Cebollín"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 70
$ws.Range("K97").Value = 6500
$ws.Range("L97").Value = 6500
$ws.Range("M97").Value = 6500
$ws.Range("N97").Value = "$/paquete 36 unidades"
$ws.Range("O97").Value = "Región Metropolitana"
$ws.Range("P97").Value = 181
$ws.Range("Q97").Value = 36
$ws.Range("R97").Value = "Hortaliza"
